$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing "Agua y Saneamiento" (and the rows below it)
# down by one
$ws.Rows(2).Insert()

# Match the formatting of the rows below (col A: bordered/bold/centered "group name"
# style, col B: plain default style) so the newly inserted row looks like every other
# data row instead of inheriting the header's look
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New group name for the inserted row
$ws.Cells.Item(2, 1).Value = "ECOLOGIA, INGENIERIA Y SOCIEDAD - EIS"

# New consultancy text block for the inserted row
$newText = @'
35. Consultoría científica: Monitoreo 13 puntos criticos para el abastecimiento de agua potable, Barbas y Cestillal, mediante la realización de aforos mensuales.  Año de inicio: 2016, Mes de inicio: 9, Año de fin: 2017, Mes de fin: 5  Idioma: Español, Ciudad: PEREIRA, Disponibilidad: Restringido, Duración: 0, Número del contrato: CONVENIO No 336 2016, Institución en la cual prestó el servicio: Corporación Autónoma Regional De Risaralda - Carder 
 39. Consultoría científica: MONITOREO DE AGUAS SUBTERRANEAS  Año de inicio: 2016, Mes de inicio: 9, Año de fin: 2017, Mes de fin: 1  Idioma: Español, Ciudad: PEREIRA, Disponibilidad: Restringido, Duración: 0, Número del contrato: CONVENIO No 336 2016, Institución en la cual prestó el servicio: Corporación Autónoma Regional De Risaralda - Carder 
 41. Consultoría científica: CONSULTORÍA NO 1302/185 DEL 2016 "ELABORACION DE ESTUDIOS DE CALIDAD Y MONITOREO DEL AGUA"  Año de inicio: 2016, Mes de inicio: 9, Año de fin: 2016, Mes de fin: 12  Idioma: Español, Ciudad: PEREIRA, Disponibilidad: Restringido, Duración: 0, Número del contrato: 1302/185 de 2016, Institución en la cual prestó el servicio: Aguas Y Aguas De Pereira 
 45. Consultoría científica: MONITOREO QUEBRADA AGUA AZUL CONVENIO ACUASEO S.A  Año de inicio: 2016, Mes de inicio: 8, Año de fin: 2017, Mes de fin: 8  Idioma: Español, Ciudad: PEREIRA, Disponibilidad: Restringido, Duración: 0, Número del contrato: 165-2016, Institución en la cual prestó el servicio: Compañía de Servicios Públicos Domiciliarios S.A E.S.P 
 46. Consultoría científica: CONVENIO 138-2015 MONITOREO QUEBRADA AGUA AZUL CONVENIO ACUASEO  Año de inicio: 2015, Mes de inicio: 8, Año de fin: 2016, Mes de fin: 8  Idioma: Español, Ciudad: PEREIRA, Disponibilidad: Restringido, Duración: 0, Número del contrato: 138-2015, Institución en la cual prestó el servicio: Compañía de Servicios Públicos Domiciliarios S.A E.S.P 
 58. Consultoría científica: CONVENIO 111 Realizar el seguimiento y monitoreo de las condiciones hidroclimatológicas en las quebradas Agua Azul, Caño NN y Quebrada La Estrella (zonas de captación superficial de agua cruda)  Año de inicio: 2018, Mes de inicio: 8, Año de fin: 2018, Mes de fin: 12  Idioma: Español, Ciudad: PEREIRA, Disponibilidad: Restringido, Duración: 0, Número del contrato: CONVENIO 111, Institución en la cual prestó el servicio: Compañía de Servicios Públicos Domiciliarios S.A E.S.P 
 68. Consultoría científica: EEP_37 de 2020_CAMPAÑA DE MONITOREO PARA EL SEGUIMIENTO DE LA CALIDAD DEL AGUA Y EL ECOSISTEMA, TRAMO DE LOCALIZACIÓN ANTES Y DESPUÉS DE BOCATOMA BELMONTE  Año de inicio: 2020, Mes de inicio: 3, Año de fin: 2020, Mes de fin: 12  Idioma: Español, Ciudad: PEREIRA, Disponibilidad: Restringido, Duración: 10, Número del contrato: No.37 de 2020, Institución en la cual prestó el servicio: EMPRESA DE ENERGIA DE PEREIRA S.A. ESP. 
 84. Consultoría científica: Monitoreo Limnológico y Calidad del agua del río Manso  Año de inicio: 2012, Mes de inicio: 1, Año de fin: 2018, Mes de fin: 12  Idioma: Español, Ciudad: NORCASIA, Disponibilidad: Restringido, Duración: 0, Número del contrato: 2014IQA, Institución en la cual prestó el servicio: INGENIEROS QUÍMICOS Y ASOCIADOS S.A.S. 
 85. Consultoría científica: Monitoreo Limnológico y Calidad del agua del río La Miel  Año de inicio: 2014, Mes de inicio: 1, Año de fin: 2018, Mes de fin: 12  Idioma: Español, Ciudad: NORCASIA, Disponibilidad: Restringido, Duración: 0, Número del contrato: 2014IQA, Institución en la cual prestó el servicio: INGENIEROS QUÍMICOS Y ASOCIADOS S.A.S. 
 86. Consultoría científica: Monitoreo Limnológico y Calidad del agua del río Guarinó  Año de inicio: 2015, Mes de inicio: 1, Año de fin: 2018, Mes de fin: 12  Idioma: Español, Ciudad: NORCASIA, Disponibilidad: Restringido, Duración: 0, Número del contrato: , Institución en la cual prestó el servicio: INGENIEROS QUÍMICOS Y ASOCIADOS S.A.S. 
 91. Consultoría científica: Monitoreo Limnológico y Calidad del agua de las quebradas que atraviesas el trasvase del río Guarinó  Año de inicio: 2016, Mes de inicio: 1, Año de fin: 2018, Mes de fin: 12  Idioma: Español, Ciudad: NORCASIA, Disponibilidad: Restringido, Duración: 0, Número del contrato: 2014IQA, Institución en la cual prestó el servicio: INGENIEROS QUÍMICOS Y ASOCIADOS S.A.S.
'@
$ws.Cells.Item(2, 2).Value = $newText
$ws.Rows(2).AutoFit()

# Remove the old last row ("Limnología y Recursos Hídricos" / its text), now pushed
# down to row 5
$ws.Rows(5).Delete()
